$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left table (Compressão): rename header and fill compression ratios ---
# C5 already holds "Compressão" (text unchanged)

# Row 6 (airport)
$ws.Range("D6").Value = "35.16:1"
$ws.Range("E6").Value = "22.88:1"
$ws.Range("F6").Value = "15.36:1"

# Row 7 (geometric)
$ws.Range("D7").Value = "87.16:1"
$ws.Range("E7").Value = "64.79:1"
$ws.Range("F7").Value = "49.47:1"

# Row 8 (nature)
$ws.Range("D8").Value = "48.04:1"
$ws.Range("E8").Value = "18.55:1"
$ws.Range("F8").Value = "14.01:1"

# --- Right table: rename header to Qualidade-Subjetiva and fill subjective quality ratings ---
$ws.Range("I5").Value = "Qualidade-Subjetiva"

# Row 6 (airport)
$ws.Range("J6").Value = "Média"
$ws.Range("K6").Value = "Média"
$ws.Range("L6").Value = "Alta"

# Row 7 (geometric)
$ws.Range("J7").Value = "Alta"
$ws.Range("K7").Value = "Alta"
$ws.Range("L7").Value = "Alta"

# Row 8 (nature)
$ws.Range("J8").Value = "Baixa"
$ws.Range("K8").Value = "Média"
$ws.Range("L8").Value = "Alta"

# New trailing blank (but formatted) cell at M6, matching the extended table border
$ws.Range("M6").Value = ""

# Apply the same centered formatting used throughout the table to all newly
# populated / touched cells so they match the existing look (style used by C5:L8)
$newCells = @("D6","E6","F6","D7","E7","F7","D8","E8","F8","J6","K6","L6","J7","K7","L7","J8","K8","L8","M6")
foreach ($addr in $newCells) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# Column I needs to widen a bit to fit the new, longer header text
$ws.Columns("I").ColumnWidth = 17.5

# Update the view: scroll so column F is at the left edge, zoom in, and
# leave the selection on M9 (matches the saved view state of the edit)
$excel.Goto($ws.Range("F1"), $true)
$excel.ActiveWindow.Zoom = 151
$ws.Range("M9").Select() | Out-Null
